# Trade #26 closed at 2026-02-17 04:08:56 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.32   # Current Capital
$summary.Range("B4").Value = 0.32      # Total P&L $
$summary.Range("B6").Value = 26        # Total Trades
$summary.Range("B7").Value = 11        # Winning Trades
$summary.Range("B9").Value = 42.31     # Win Rate %

# --- Strategy Status sheet (MarketMaking row, row 4) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.32
$status.Range("D4").Value = 26
$status.Range("E4").Value = 0.32
$status.Range("F4").Value = 0.32
$status.Range("G4").Value = 42.31

# --- New trade row appended to "All Trades" and "MarketMaking" sheets ---
# Note: column B holds a date formatted as plain text (e.g. "2026-02-17"),
# like all the rows above it. A bare Value assignment of an ISO-formatted
# date string gets auto-recognized and converted to a real date by Excel's
# smart input, so it is entered with a leading apostrophe to force it to be
# stored as literal text (the normal Excel way of entering "text that looks
# like a date"), matching the existing rows.
$newRow = @(26, "'2026-02-17", "04:08:50", "MarketMaking", "UP", 0.57, 0.58, "CLOSED", 1.7544, 0.01, 100.32, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.12)

$allTrades = $wb.Worksheets.Item("All Trades")
$marketMaking = $wb.Worksheets.Item("MarketMaking")

foreach ($ws in @($allTrades, $marketMaking)) {
    $row = 27
    for ($col = 1; $col -le $newRow.Length; $col++) {
        $ws.Cells.Item($row, $col).Value = $newRow[$col - 1]
    }
}
